$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) is treated as text so values like "28.850.28"
# or trailing-zero decimals are not reinterpreted as numbers/dates.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '28.850.28'
$ws.Range("E2").Value = '  +8.04%  '
$ws.Range("D3").Value = '1.810.08'
$ws.Range("E3").Value = '  +5.06%  '
$ws.Range("D4").Value = '0.9994'
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("D5").Value = '247.28'
$ws.Range("E5").Value = '  +3.08%  '
$ws.Range("D6").Value = '0.9994'
$ws.Range("E6").Value = '  +0.09%  '
$ws.Range("D7").Value = '0.4971'
$ws.Range("E7").Value = '  +2.85%  '
$ws.Range("D8").Value = '0.2787'
$ws.Range("E8").Value = '  +8.28%  '
$ws.Range("D9").Value = '0.06416'
$ws.Range("E9").Value = '  +3.93%  '
$ws.Range("D10").Value = '1.814.20'
$ws.Range("E10").Value = '  +5.37%  '
$ws.Range("D11").Value = '16.81'
$ws.Range("E11").Value = '  +6.00%  '
$ws.Range("D12").Value = '0.07082'
$ws.Range("E12").Value = '  +3.70%  '
$ws.Range("D13").Value = '0.6469'
$ws.Range("E13").Value = '  +7.16%  '
$ws.Range("D14").Value = '84.11'
$ws.Range("E14").Value = '  +9.45%  '
$ws.Range("D15").Value = '4.699'
$ws.Range("E15").Value = '  +5.52%  '
$ws.Range("D16").Value = '28.866.12'
$ws.Range("E16").Value = '  +8.79%  '
$ws.Range("D17").Value = '0.9996'
$ws.Range("E17").Value = '  +0.12%  '
$ws.Range("D18").Value = '0.000007337'
$ws.Range("E18").Value = '  +2.77%  '
$ws.Range("D19").Value = '0.9992'
$ws.Range("E19").Value = '  +0.09%  '
$ws.Range("D20").Value = '12.28'
$ws.Range("E20").Value = '  +8.15%  '
$ws.Range("D21").Value = '2.048.45'
$ws.Range("E21").Value = '  +5.41%  '
$ws.Range("D22").Value = '4.590'
$ws.Range("E22").Value = '  +4.05%  '
$ws.Range("D23").Value = '8.878'
$ws.Range("E23").Value = '  +3.76%  '
$ws.Range("D24").Value = '5.346'
$ws.Range("E24").Value = '  +5.87%  '
$ws.Range("D25").Value = '142.30'
$ws.Range("E25").Value = '  +2.22%  '
$ws.Range("D26").Value = '129.50'
$ws.Range("E26").Value = '  +21.77%  '
$ws.Range("D27").Value = '16.40'
$ws.Range("E27").Value = '  +7.89%  '
$ws.Range("D28").Value = '1.887'
$ws.Range("E28").Value = '  +6.79%  '
$ws.Range("D29").Value = '1.410'
$ws.Range("E29").Value = '  +3.13%  '
$ws.Range("D30").Value = '4.141'
$ws.Range("E30").Value = '  +3.47%  '
$ws.Range("D31").Value = '0.08366'
$ws.Range("E31").Value = '  +6.01%  '
$ws.Range("D32").Value = '3.811'
$ws.Range("E32").Value = '  +4.44%  '
$ws.Range("D33").Value = '0.04965'
$ws.Range("E33").Value = '  +10.70%  '
$ws.Range("D34").Value = '1.096'
$ws.Range("E34").Value = '  +9.86%  '
$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D35").Value = '0.6737'
$ws.Range("E35").Value = '  +9.44%  '
$ws.Range("B36").Value = 'HuobiToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D36").Value = '2.688'
$ws.Range("E36").Value = '  +3.66%  '
$ws.Range("D37").Value = '2.310'
$ws.Range("E37").Value = '  +15.25%  '
$ws.Range("D38").Value = '2.750'
$ws.Range("E38").Value = '  +12.50%  '
$ws.Range("D39").Value = '0.9538'
$ws.Range("E39").Value = '  +1.93%  '
$ws.Range("D40").Value = '6.138'
$ws.Range("E40").Value = '  +9.45%  '
$ws.Range("D41").Value = '0.01593'
$ws.Range("E41").Value = '  +7.03%  '
$ws.Range("D42").Value = '0.9995'
$ws.Range("E42").Value = '  +0.13%  '
$ws.Range("D43").Value = '0.4096'
$ws.Range("E43").Value = '  +7.15%  '
$ws.Range("D44").Value = '100.05'
$ws.Range("E44").Value = '  +0.19%  '
$ws.Range("D45").Value = '7.165'
$ws.Range("E45").Value = '  +5.83%  '
$ws.Range("E46").Value = '  +5.98%  '
$ws.Range("D47").Value = '0.05512'
$ws.Range("E47").Value = '  +2.79%  '
$ws.Range("B48").Value = 'Elrond'
$ws.Range("C48").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D48").Value = '31.69'
$ws.Range("E48").Value = '  +5.75%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = '8.097'
$ws.Range("E49").Value = '  +2.43%  '
$ws.Range("B50").Value = 'Decentraland'
$ws.Range("C50").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D50").Value = '0.3640'
$ws.Range("E50").Value = '  +9.05%  '
$ws.Range("B51").Value = 'NEARProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D51").Value = '1.313'
$ws.Range("E51").Value = '  +6.02%  '
